$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a cell to be written as literal text (preserve formatting such as
# "1.00" or multi-dot numbers) regardless of Excel's automatic type detection,
# then restore the default "Normal" style so no stray formatting is introduced.
function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "66.646.71"
Set-TextCell "E2" "  +3.10%  "
Set-TextCell "D3" "3.195.34"
Set-TextCell "E3" "  +1.51%  "
Set-TextCell "E4" "  +0.04%  "
Set-TextCell "D5" "597.67"
Set-TextCell "E5" "  +3.52%  "
Set-TextCell "D6" "155.19"
Set-TextCell "E6" "  +4.25%  "
Set-TextCell "D7" "1.00"
Set-TextCell "E7" "  -0.01%  "
Set-TextCell "D8" "0.557"
Set-TextCell "E8" "  +5.95%  "
Set-TextCell "D9" "3.182.75"
Set-TextCell "E9" "  +1.17%  "
Set-TextCell "E10" "  +1.55%  "
Set-TextCell "D11" "5.95"
Set-TextCell "E11" "  -2.76%  "
Set-TextCell "E12" "  +3.53%  "
Set-TextCell "E13" "  +3.29%  "
Set-TextCell "D14" "39.41"
Set-TextCell "E14" "  +6.03%  "
Set-TextCell "D15" "3.718.85"
Set-TextCell "E15" "  +1.45%  "
Set-TextCell "B16" "Polkadot"
Set-TextCell "C16" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D16" "7.51"
Set-TextCell "E16" "  +5.29%  "
Set-TextCell "B17" "WrappedBTC"
Set-TextCell "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D17" "66.582.08"
Set-TextCell "E17" "  +2.89%  "
Set-TextCell "D18" "3.198.98"
Set-TextCell "E18" "  +2.08%  "
Set-TextCell "E19" "  +0.50%  "
Set-TextCell "D20" "518.65"
Set-TextCell "E20" "  +2.67%  "
Set-TextCell "D21" "15.41"
Set-TextCell "E21" "  +3.69%  "
Set-TextCell "D22" "0.741"
Set-TextCell "E22" "  +3.85%  "
Set-TextCell "D23" "8.10"
Set-TextCell "E23" "  +5.13%  "
Set-TextCell "D24" "15.01"
Set-TextCell "E24" "  -1.49%  "
Set-TextCell "E25" "  +2.11%  "
Set-TextCell "E26" "  +0.05%  "
Set-TextCell "E27" "  +4.82%  "
Set-TextCell "E28" "  +4.22%  "
Set-TextCell "E29" "  +8.91%  "
Set-TextCell "D30" "7.15"
Set-TextCell "E30" "  +15.61%  "
Set-TextCell "E31" "  +4.89%  "
Set-TextCell "D32" "28.40"
Set-TextCell "E32" "  +3.12%  "
Set-TextCell "D33" "1.24"
Set-TextCell "E33" "  +3.67%  "
Set-TextCell "D34" "1.00"
Set-TextCell "E34" "  +0.17%  "
Set-TextCell "E35" "  +1.48%  "
Set-TextCell "D36" "517.71"
Set-TextCell "E36" "  +7.55%  "
Set-TextCell "D37" "55.04"
Set-TextCell "E37" "  +0.90%  "
Set-TextCell "D38" "0.0905"
Set-TextCell "E38" "  +1.51%  "
Set-TextCell "D39" "0.0426"
Set-TextCell "E39" "  +2.66%  "
Set-TextCell "D40" "0.130"
Set-TextCell "E40" "  +13.13%  "
Set-TextCell "D41" "8.94"
Set-TextCell "E41" "  +2.69%  "
Set-TextCell "D42" "2.94"
Set-TextCell "E42" "  +0.30%  "
Set-TextCell "E43" "  +8.06%  "
Set-TextCell "D44" "0.0₃0677"
Set-TextCell "E44" "  +16.37%  "
Set-TextCell "D45" "2.45"
Set-TextCell "E45" "  +1.58%  "
Set-TextCell "D46" "2.914.13"
Set-TextCell "E46" "  -2.86%  "
Set-TextCell "D47" "28.86"
Set-TextCell "E48" "  +3.61%  "
Set-TextCell "B49" "ThetaToken"
Set-TextCell "C49" "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell "D49" "2.36"
Set-TextCell "E49" "  +5.99%  "
Set-TextCell "B50" "USDe"
Set-TextCell "C50" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell "D50" "0.999"
Set-TextCell "E50" "  -0.01%  "
Set-TextCell "D51" "2.66"
Set-TextCell "E51" "  +9.42%  "
